# Rerouting.xlsx sample-test workbook was repurposed from a "One Way Anova"
# example into the "Reroute To Sea" customer test, so the (only) worksheet
# just needs renaming to match the new scenario it now documents.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("One Way Anova")
$ws.Name = "Reroute To Sea"
